# Scheduled runner: refresh market-price-derived profit figures across the
# per-job "Profits" sheets (currentAveragePrice* / LevePrice* / LeveProfit*
# columns, H..N) to their latest observed values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2481.5454
$ws.Range("I19").Value = 1333.6666
$ws.Range("K19").Value = 1333.6666
$ws.Range("M19").Value = -1158.6666

$ws.Range("H40").Value = 3984.8
$ws.Range("I40").Value = 3678.5715
$ws.Range("K40").Value = 3678.5715
$ws.Range("M40").Value = -3503.5715

$ws.Range("H53").Value = 1070.6428
$ws.Range("I53").Value = 510.66666
$ws.Range("J53").Value = 2078.6
$ws.Range("K53").Value = 510.66666
$ws.Range("L53").Value = 2078.6
$ws.Range("M53").Value = 126.33334
$ws.Range("N53").Value = -3352.6

$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()

$ws.Range("H55").Value = 740
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()

$ws.Range("H76").Value = 8291.4
$ws.Range("J76").Value = 6653.222
$ws.Range("L76").Value = 6653.222
$ws.Range("N76").Value = -7283.222

$ws.Range("H79").Value = 8291.4
$ws.Range("J79").Value = 6653.222
$ws.Range("L79").Value = 6653.222
$ws.Range("N79").Value = -8837.222

$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("N97").ClearContents()

$ws.Range("H100").Value = 5246.1665
$ws.Range("I100").Value = 1177.6364
$ws.Range("J100").Value = 50000
$ws.Range("K100").Value = 1177.6364
$ws.Range("L100").Value = 50000
$ws.Range("M100").Value = -636.6364000000001
$ws.Range("N100").Value = -51082

$ws.Range("H112").Value = 1995.5
$ws.Range("I112").Value = 1006.3333
$ws.Range("J112").Value = 2223.7693
$ws.Range("K112").Value = 3018.9999
$ws.Range("L112").Value = 6671.3079
$ws.Range("M112").Value = -1910.9999
$ws.Range("N112").Value = -8887.3079

$ws.Range("H115").Value = 491.77777
$ws.Range("I115").Value = 491.77777
$ws.Range("K115").Value = 1475.33331
$ws.Range("M115").Value = 91.66669000000002

$ws.Range("H116").Value = 4941.8
$ws.Range("J116").Value = 5861
$ws.Range("L116").Value = 5861
$ws.Range("N116").Value = -12745

$ws.Range("H132").Value = 973.64703
$ws.Range("I132").Value = 851.6667
$ws.Range("K132").Value = 2555.0001
$ws.Range("M132").Value = -25.0001000000002

$ws.Range("H138").Value = 2479.9744
$ws.Range("I138").Value = 1440.3334
$ws.Range("J138").Value = 2727.508
$ws.Range("K138").Value = 4321.0002
$ws.Range("L138").Value = 8182.523999999999
$ws.Range("M138").Value = 818.9997999999996
$ws.Range("N138").Value = -18462.524

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9096867
$ws.Range("I32").Value = 12502804
$ws.Range("K32").Value = 12502804
$ws.Range("M32").Value = -12502517

$ws.Range("H53").Value = 5000
$ws.Range("J53").Value = 5000
$ws.Range("L53").Value = 5000
$ws.Range("N53").Value = -6364

$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()

$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents()

$ws.Range("H63").Value = 4949.7896
$ws.Range("I63").Value = 2366.818
$ws.Range("K63").Value = 2366.818
$ws.Range("M63").Value = -1680.818

$ws.Range("H66").Value = 4949.7896
$ws.Range("I66").Value = 2366.818
$ws.Range("K66").Value = 11834.09
$ws.Range("M66").Value = -8402.09

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H49").Value = 5489.5
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 5489.5
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 5489.5
$ws.Range("M49").ClearContents()
$ws.Range("N49").Value = -5967.5

$ws.Range("H54").Value = 19699.5
$ws.Range("J54").Value = 8799
$ws.Range("L54").Value = 8799
$ws.Range("N54").Value = -9767

$ws.Range("H107").Value = 2923.9443
$ws.Range("I107").Value = 2506.1333
$ws.Range("K107").Value = 2506.1333
$ws.Range("M107").Value = -586.1333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 561291.25
$ws.Range("I31").Value = 9131.362999999999
$ws.Range("J31").Value = 1168667.1
$ws.Range("K31").Value = 9131.362999999999
$ws.Range("L31").Value = 1168667.1
$ws.Range("M31").Value = -8836.362999999999
$ws.Range("N31").Value = -1169257.1

$ws.Range("H34").Value = 561291.25
$ws.Range("I34").Value = 9131.362999999999
$ws.Range("J34").Value = 1168667.1
$ws.Range("K34").Value = 9131.362999999999
$ws.Range("L34").Value = 1168667.1
$ws.Range("M34").Value = -8929.362999999999
$ws.Range("N34").Value = -1169071.1

$ws.Range("H39").Value = 9333
$ws.Range("I39").Value = 9333
$ws.Range("K39").Value = 9333
$ws.Range("M39").Value = -8942

$ws.Range("H49").Value = 9333
$ws.Range("I49").Value = 9333
$ws.Range("K49").Value = 9333
$ws.Range("M49").Value = -9151

$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("N55").ClearContents()

$ws.Range("H62").Value = 1433556.9
$ws.Range("I62").Value = 2005879.6
$ws.Range("K62").Value = 2005879.6
$ws.Range("M62").Value = -2005255.6

$ws.Range("H65").Value = 1433556.9
$ws.Range("I65").Value = 2005879.6
$ws.Range("K65").Value = 10029398
$ws.Range("M65").Value = -10026278

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 61495
$ws.Range("J37").Value = 61495
$ws.Range("L37").Value = 184485
$ws.Range("N37").Value = -184709

$ws.Range("H54").Value = 22857.143
$ws.Range("J54").Value = 24230.77
$ws.Range("L54").Value = 72692.31
$ws.Range("N54").Value = -73810.31

$ws.Range("H55").Value = 2874.6667
$ws.Range("I55").Value = 4102
$ws.Range("J55").Value = 420
$ws.Range("K55").Value = 12306
$ws.Range("L55").Value = 1260
$ws.Range("M55").Value = -12129
$ws.Range("N55").Value = -1614

$ws.Range("H116").Value = 800
$ws.Range("I116").Value = 800
$ws.Range("K116").Value = 2400
$ws.Range("M116").Value = 1042

$ws.Range("H131").Value = 21310.953
$ws.Range("J131").Value = 19876.5
$ws.Range("L131").Value = 59629.5
$ws.Range("N131").Value = -69709.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 22015.5
$ws.Range("J47").Value = 22015.5
$ws.Range("L47").Value = 22015.5
$ws.Range("N47").Value = -23151.5

$ws.Range("H49").Value = 38498.332
$ws.Range("J49").Value = 38498.332
$ws.Range("L49").Value = 38498.332
$ws.Range("N49").Value = -38866.332

$ws.Range("H70").Value = 27799.2
$ws.Range("I70").Value = 27799.2
$ws.Range("K70").Value = 27799.2
$ws.Range("M70").Value = -27529.2

$ws.Range("H73").Value = 27799.2
$ws.Range("I73").Value = 27799.2
$ws.Range("K73").Value = 27799.2
$ws.Range("M73").Value = -26863.2

$ws.Range("H80").Value = 4282.6
$ws.Range("I80").Value = 3647.4443
$ws.Range("J80").Value = 9999
$ws.Range("K80").Value = 3647.4443
$ws.Range("L80").Value = 9999
$ws.Range("M80").Value = -2649.4443
$ws.Range("N80").Value = -11995

$ws.Range("H83").Value = 4282.6
$ws.Range("I83").Value = 3647.4443
$ws.Range("J83").Value = 9999
$ws.Range("K83").Value = 18237.2215
$ws.Range("L83").Value = 49995
$ws.Range("M83").Value = -13245.2215
$ws.Range("N83").Value = -59979

$ws.Range("H102").Value = 2414.3157
$ws.Range("I102").Value = 1878.4073
$ws.Range("J102").Value = 3729.7273
$ws.Range("K102").Value = 1878.4073
$ws.Range("L102").Value = 3729.7273
$ws.Range("M102").Value = -256.4073000000001
$ws.Range("N102").Value = -6973.7273

$ws.Range("H123").Value = 54000
$ws.Range("J123").Value = 54000
$ws.Range("L123").Value = 54000
$ws.Range("N123").Value = -58900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1758.4445
$ws.Range("J22").Value = 1250.1818
$ws.Range("L22").Value = 1250.1818
$ws.Range("N22").Value = -1840.1818

$ws.Range("H27").Value = 1758.4445
$ws.Range("J27").Value = 1250.1818
$ws.Range("L27").Value = 1250.1818
$ws.Range("N27").Value = -1464.1818

$ws.Range("H46").Value = 3338.8462
$ws.Range("I46").Value = 2080
$ws.Range("K46").Value = 2080
$ws.Range("M46").Value = -1892

$ws.Range("H53").Value = 41900
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents()

$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()

$ws.Range("H55").Value = 58824320
$ws.Range("I55").Value = 90909820
$ws.Range("J55").Value = 893.3333
$ws.Range("K55").Value = 90909820
$ws.Range("L55").Value = 893.3333
$ws.Range("M55").Value = -90909647
$ws.Range("N55").Value = -1239.3333
